$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.01

$ws.Range("F3").Value = 0.49
$ws.Range("G3").Value = 0.51

$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

$ws.Range("G5").Value = 0.98
$ws.Range("H5").Value = 0.01

$ws.Range("G6").Value = 0.01
